$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.941.94"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.34%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.856.26"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.87%  "

# Row 4
$ws.Range("E4").Value = "  -0.15%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "471.55"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +10.13%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.62"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +10.07%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.630"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.69%  "

# Row 8
$ws.Range("E8").Value = "  -0.12%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.744"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.52%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.154"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.66%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0000311"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -7.64%  "

# Row 12
$ws.Range("E12").Value = "  +4.92%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.38"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.47%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.488.13"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.07%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.78"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -5.02%  "

# Row 16
$ws.Range("B16").Value = "TRON"
$ws.Range("C16").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.137"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.37%  "

# Row 17
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.820.16"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.49%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "20.10"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.02%  "

# Row 19
$ws.Range("E19").Value = "  +4.67%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "67.164.22"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.25%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "431.50"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.57%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "14.88"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.60%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.33"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +6.37%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "88.40"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.47%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.60"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +9.03%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "37.98"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.88%  "

# Row 27
$ws.Range("E27").Value = "  +5.95%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.91"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.96%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.55"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.05%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "730.72"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.53%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "13.85"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.85%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.134"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +6.95%  "

# Row 33
$ws.Range("E33").Value = "  +0.73%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "43.29"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +11.14%  "

# Row 35
$ws.Range("E35").Value = "  +7.09%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "58.11"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.40%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.999"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.09%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.46"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -6.15%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0484"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.82%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.347"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +7.64%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.91"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.04%  "

# Row 42
$ws.Range("E42").Value = "  +2.96%  "

# Row 43
$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.00"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.10%  "

# Row 44
$ws.Range("B44").Value = "PEPE"
$ws.Range("C44").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0₃0672"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -7.23%  "

# Row 45
$ws.Range("B45").Value = "Fetch.AI"
$ws.Range("C45").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.55"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +5.71%  "

# Row 46
$ws.Range("E46").Value = "  +1.89%  "

# Row 47
$ws.Range("E47").Value = "  +5.14%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.15"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.81%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.21"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.05%  "

# Row 50
$ws.Range("E50").Value = "  +1.54%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "143.46"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.49%  "
